# Map of python files.pptx - apply diagram edits
# - Reposition the "dbconfig" and "dbconnection" boxes
# - Rename/reposition "UserConnection" box to "User"
# - Reroute / resize three connector arrows (dbconnection->Population,
#   dbconnection->User, Population->main)
# - Add a new connector arrow from "dbconfig" down to "dbconnection"
# - Add a new "Population" box plus its feeding arrow

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id=5 "TextBox 4" (dbconfig) : move up ---
$s.Shapes.Item(2).Top = 31.10212698425197

# --- Shape id=6 "TextBox 5" (dbconnection) : move up ---
$s.Shapes.Item(3).Top = 141.09023622047243

# --- Shape id=7 "TextBox 6" (UserConnection -> User) : reposition + retext ---
$shp7 = $s.Shapes.Item(4)
$shp7.TextFrame.TextRange.Delete()
$shp7.TextFrame.TextRange.Text = "User"
$shp7.Left = 400.65300012598425
$shp7.Top = 258.22984351968506
$shp7.Width = 158.693937007874
$shp7.Height = 32.17503937007874

# --- Shape id=17 "Straight Arrow Connector 16" : reroute ---
$shp17 = $s.Shapes.Item(7)
$shp17.Left = 220.65276390551182
$shp17.Top = 187.3963859527559
$shp17.Width = 0.0
$shp17.Height = 55.00299312598425

# --- Shape id=18 "Straight Arrow Connector 17" : reroute ---
$shp18 = $s.Shapes.Item(8)
$shp18.Left = 310.5303937007874
$shp18.Top = 158.67937507874015
$shp18.Width = 153.79614173228347
$shp18.Height = 90.17086614173229

# --- Shape id=28 "Straight Arrow Connector 27" : reroute ---
$shp28 = $s.Shapes.Item(12)
$shp28.Left = 220.65267716535433
$shp28.Top = 296.42
$shp28.Width = 166.04110236220473
$shp28.Height = 106.19228346456693

# --- New connector: "Straight Arrow Connector 13" (main -> dbconfig) ---
# Duplicate an existing styled connector so the line/arrow formatting matches.
$newConnRange = $shp17.Duplicate()
$newConn = $newConnRange.Item(1)
$newConn.Name = "Straight Arrow Connector 13"
$newConn.Left = 220.89779527559054
$newConn.Top = 77.1428346456693
$newConn.Width = 0.0
$newConn.Height = 51.42858367716535

# --- New shape: "TextBox 18" (Population) ---
# Duplicate an existing styled textbox so the shape/text formatting matches.
$newTbRange = $s.Shapes.Item(3).Duplicate()
$newTb = $newTbRange.Item(1)
$newTb.Name = "TextBox 18"
$newTb.TextFrame.TextRange.Delete()
$newTb.TextFrame.TextRange.Text = "Population"
$newTb.Left = 141.30574803149605
$newTb.Top = 256.2649696299212
$newTb.Width = 158.693937007874
$newTb.Height = 32.17503937007874
